# Keep the COM layer from stamping synthetic rsids where it can be
# helped, and always re-fetch Range objects from the live Paragraphs
# collection right before using them (stale Range handles captured
# right after an insert make the shim mis-attribute bookkeeping to the
# wrong paragraph).
$word.Options.StoreRSIDOnSave = $false

$d = $word.ActiveDocument

# The document ends with a SourceCode paragraph whose last line of
# verbatim console output is "## 6 2007-06-01       31     6 Summer".
# Append a brand-new SourceCode paragraph after it with two more
# `library(...)` calls, matching the formatting already used for the
# `library(readr)` / `library(tidyverse)` paragraph at the top of the
# document: FunctionTok style on the function name, NormalTok on the
# parenthesized argument, and a manual line break between the two
# statements.

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$tail = $lastPara.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

# The freshly inserted paragraph already inherited pStyle "SourceCode"
# from its predecessor, which is exactly what we want.
$p1 = $d.Paragraphs.Count
$r = $d.Paragraphs($p1).Range
$r.Collapse(0)
$r.InsertAfter("library(ggplot2)")

# Style "library" -> FunctionTok
$fr = $d.Paragraphs($p1).Range
$fr.Find.ClearFormatting()
$fr.Find.Replacement.ClearFormatting()
$fr.Find.Replacement.Style = "FunctionTok"
$fr.Find.Execute("library", $true, $false, $false, $false, $false, $true, 1, $false, "library", 2)

# Style "(ggplot2)" -> NormalTok
$fr2 = $d.Paragraphs($p1).Range
$fr2.Find.ClearFormatting()
$fr2.Find.Replacement.ClearFormatting()
$fr2.Find.Replacement.Style = "NormalTok"
$fr2.Find.Execute("(ggplot2)", $true, $false, $false, $false, $false, $true, 1, $false, "(ggplot2)", 2)

# Manual line break (<w:br/>), then the second library() call in the
# same paragraph.
$br = $d.Paragraphs($d.Paragraphs.Count).Range
$br.Collapse(0)
$br.InsertAfter([char]11)

$p2 = $d.Paragraphs.Count
$r2 = $d.Paragraphs($p2).Range
$r2.Collapse(0)
$r2.InsertAfter("library(ggthemes)")

# Style "library" -> FunctionTok
$fr3 = $d.Paragraphs($p2).Range
$fr3.Find.ClearFormatting()
$fr3.Find.Replacement.ClearFormatting()
$fr3.Find.Replacement.Style = "FunctionTok"
$fr3.Find.Execute("library", $true, $false, $false, $false, $false, $true, 1, $false, "library", 2)

# Style "(ggthemes)" -> NormalTok
$fr4 = $d.Paragraphs($p2).Range
$fr4.Find.ClearFormatting()
$fr4.Find.Replacement.ClearFormatting()
$fr4.Find.Replacement.Style = "NormalTok"
$fr4.Find.Execute("(ggthemes)", $true, $false, $false, $false, $false, $true, 1, $false, "(ggthemes)", 2)

Write-Output "Done. Paragraphs=$($d.Paragraphs.Count)"
